$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking values (price/volume figures) must be forced to Text format
# so Excel stores the exact original string instead of auto-converting it to a
# Double/Percent number (which would lose formatting/precision).
$numericCells = @('D2', 'E2', 'D3', 'E3', 'D4', 'E4', 'D5', 'E5', 'E6', 'D7', 'E7', 'D8', 'E8', 'D9', 'E9', 'D10', 'E10', 'D11', 'E11', 'D12', 'E12', 'D13', 'E13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'D17', 'E17', 'D18', 'E18', 'D19', 'E19', 'D20', 'E20', 'D21', 'E21', 'D22', 'E22', 'D23', 'E23', 'D24', 'E24', 'D25', 'E25', 'E26', 'D27', 'E27', 'E39', 'D40', 'E40', 'D41', 'E41', 'D42', 'E42', 'D43', 'E43', 'D44', 'E44', 'D45', 'E45', 'D46', 'E46', 'E47', 'D48', 'E48', 'E49', 'E50', 'E51')
foreach ($addr in $numericCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '310.34'
$ws.Range("E2").Value = '1.45%'
# Row 3
$ws.Range("D3").Value = '35.60'
$ws.Range("E3").Value = '-1.51%'
# Row 4
$ws.Range("D4").Value = '5.113'
$ws.Range("E4").Value = '1.33%'
# Row 5
$ws.Range("D5").Value = '0.08213'
$ws.Range("E5").Value = '4.20%'
# Row 6
$ws.Range("E6").Value = '-8.13%'
# Row 7
$ws.Range("D7").Value = '7.968'
$ws.Range("E7").Value = '-0.35%'
# Row 8
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").Value = '2.896'
$ws.Range("E8").Value = '8.87%'
# Row 9
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '0.9277'
$ws.Range("E9").Value = '-0.04%'
# Row 10
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '0.1084'
$ws.Range("E10").Value = '10.63%'
# Row 11
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '0.1921'
$ws.Range("E11").Value = '2.55%'
# Row 12
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '0.09592'
$ws.Range("E12").Value = '6.97%'
# Row 13
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '0.03572'
$ws.Range("E13").Value = '-4.94%'
# Row 14
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '0.09907'
$ws.Range("E14").Value = '-0.21%'
# Row 15
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '0.001446'
$ws.Range("E15").Value = '0.50%'
# Row 16
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '0.005734'
$ws.Range("E16").Value = '1.13%'
# Row 17
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '3.475'
$ws.Range("E17").Value = '0.41%'
# Row 18
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").Value = '4.130'
$ws.Range("E18").Value = '-0.63%'
# Row 19
$ws.Range("D19").Value = '0.3426'
$ws.Range("E19").Value = '1.75%'
# Row 20
$ws.Range("D20").Value = '0.1310'
$ws.Range("E20").Value = '-0.76%'
# Row 21
$ws.Range("D21").Value = '5.097'
$ws.Range("E21").Value = '0.30%'
# Row 22
$ws.Range("D22").Value = '0.2192'
$ws.Range("E22").Value = '-2.68%'
# Row 23
$ws.Range("D23").Value = '0.04546'
$ws.Range("E23").Value = '-0.56%'
# Row 24
$ws.Range("D24").Value = '0.001226'
$ws.Range("E24").Value = '-0.67%'
# Row 25
$ws.Range("D25").Value = '0.004797'
$ws.Range("E25").Value = '0.54%'
# Row 26
$ws.Range("E26").Value = '-3.88%'
# Row 27
$ws.Range("D27").Value = '0.0004452'
$ws.Range("E27").Value = '-6.05%'
# Row 39
$ws.Range("E39").Value = '1.97%'
# Row 40
$ws.Range("D40").Value = '0.04912'
$ws.Range("E40").Value = '-0.50%'
# Row 41
$ws.Range("D41").Value = '0.007648'
$ws.Range("E41").Value = '-1.72%'
# Row 42
$ws.Range("D42").Value = '0.009859'
$ws.Range("E42").Value = '26.15%'
# Row 43
$ws.Range("D43").Value = '0.1384'
$ws.Range("E43").Value = '-0.56%'
# Row 44
$ws.Range("D44").Value = '0.002117'
$ws.Range("E44").Value = '-1.22%'
# Row 45
$ws.Range("D45").Value = '0.01157'
$ws.Range("E45").Value = '1.18%'
# Row 46
$ws.Range("D46").Value = '0.00006528'
$ws.Range("E46").Value = '6.20%'
# Row 47
$ws.Range("E47").Value = '-0.11%'
# Row 48
$ws.Range("D48").Value = '62.27'
$ws.Range("E48").Value = '20.29%'
# Row 49
$ws.Range("E49").Value = '-16.89%'
# Row 50
$ws.Range("E50").Value = '-0.11%'
# Row 51
$ws.Range("E51").Value = '-0.11%'
